# The deck's applied colour theme is switched from the custom "Integral"
# palette over to the plain "Office Theme" palette (same effect as picking
# a different theme colour set from the Design tab in PowerPoint).
#
# PowerPoint exposes the twelve theme colour slots (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) through ThemeColorScheme.Colors(i).RGB, so
# the swap is performed one slot at a time.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeThemeColors[$i - 1]
}
